$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lines = @(
  'MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)',
  'OPTIONAL MATCH (samp:sample)-->(c)',
  'OPTIONAL MATCH (diag:diagnosis)-->(c)',
  'OPTIONAL MATCH (f:file)-[*]->(c)',
  'OPTIONAL MATCH (sf:file)-->(s)',
  'WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p',
  "WHERE demo.breed IN ['Labrador Retriever']",
  'RETURN  ',
  '    count(distinct p) AS Programs,',
  '    count(distinct s) AS Studies,',
  '    count(distinct c) AS Cases,',
  '    count(distinct samp) AS Samples,',
  '    count(distinct f) AS `Case Files`,',
  '    count(distinct sf) AS `Study Files`'
)
$newStatQuery = [string]::Join([char]10, $lines)

$ws.Cells.Item(2,3).Value2 = $newStatQuery
$ws.Cells.Item(3,3).Value2 = $newStatQuery
$ws.Cells.Item(4,3).Value2 = $newStatQuery

$ws.Rows.Item(3).RowHeight = 230.4

$ws.Cells.Item(4,2).Select()
$excel.ActiveWindow.Zoom = 100
